$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rule row 11 ("R40") in the decision table is renamed to "1"
# Force text storage (shared string), not a numeric value.
$cell = $ws.Range("B11")
$cell.NumberFormat = "@"
$cell.Value = "1"
